$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column F
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# time_taken values for each data row (2-21)
$timeTaken = @(
    "2021-10-05 10:52:53.828550",
    "2021-10-05 10:52:53.828563",
    "2021-10-05 10:52:53.828567",
    "2021-10-05 10:52:53.828570",
    "2021-10-05 10:52:53.828574",
    "2021-10-05 10:52:53.828577",
    "2021-10-05 10:52:53.828580",
    "2021-10-05 10:52:53.828583",
    "2021-10-05 10:52:53.828587",
    "2021-10-05 10:52:53.828590",
    "2021-10-05 10:52:53.828593",
    "2021-10-05 10:52:53.828596",
    "2021-10-05 10:52:53.828599",
    "2021-10-05 10:52:53.828602",
    "2021-10-05 10:52:53.828605",
    "2021-10-05 10:52:53.828608",
    "2021-10-05 10:52:53.828612",
    "2021-10-05 10:52:53.828615",
    "2021-10-05 10:52:53.828618",
    "2021-10-05 10:52:53.828621"
)

for ($i = 0; $i -lt $timeTaken.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timeTaken[$i]
}
